$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lookups")

# Insert two new columns before the existing "question_code" pair (H:I) to make
# room for a new "GenderB" lookup pair; this shifts question_code -> J:K and
# result_type -> L:M.
$ws.Range("H1:I1").EntireColumn.Insert()

# Populate the new GenderB header + male/female lookup values (mirrors the
# existing Gender pair in columns F:G / H:I).
$ws.Range("H1").Value = "GenderB"
$ws.Range("H2").Value = "male"
$ws.Range("I2").Value = 1
$ws.Range("H3").Value = "female"
$ws.Range("I3").Value = 2

# Extend the named range that spans the lookup header row to include the
# newly added columns.
$wb.Names.Item("cuts_head").RefersTo = "='Lookups'!`$F`$1:`$M`$1"
